$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Rtn4"
$ws.Range("C2").Value = "Rtn4rl1"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 35.71508866666667
$ws.Range("H2").Value = 107.145266
$ws.Range("I2").Value = 0.1390302752364672
$ws.Range("J2").Value = 0.1390302752364672
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.620945
$ws.Range("N2").Value = 4.862835
$ws.Range("O2").Value = 0.4497939646201873
$ws.Range("P2").Value = 0.4497939646201873
$ws.Range("Q2").Value = 57.89219439879
$ws.Range("R2").Value = 521.0297495891101
$ws.Range("S2").Value = 0.06253497870084641
$ws.Range("T2").Value = 0.06253497870084643

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Rtn4"
$ws.Range("C3").Value = "Rtn4rl1"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 35.71508866666667
$ws.Range("H3").Value = 107.145266
$ws.Range("I3").Value = 0.1390302752364672
$ws.Range("J3").Value = 0.1390302752364672
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.421245
$ws.Range("N3").Value = 4.263735
$ws.Range("O3").Value = 0.3943794658342005
$ws.Range("P3").Value = 0.3943794658342005
$ws.Range("Q3").Value = 50.75989119205667
$ws.Range("R3").Value = 456.83902072851
$ws.Range("S3").Value = 0.05483068568253979
$ws.Range("T3").Value = 0.0548306856825398

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Rtn4"
$ws.Range("C4").Value = "Rtn4rl1"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 35.71508866666667
$ws.Range("H4").Value = 107.145266
$ws.Range("I4").Value = 0.1390302752364672
$ws.Range("J4").Value = 0.1390302752364672
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.5615599999999999
$ws.Range("N4").Value = 1.68468
$ws.Range("O4").Value = 0.1558265695456122
$ws.Range("P4").Value = 0.1558265695456122
$ws.Range("Q4").Value = 20.05616519165334
$ws.Range("R4").Value = 180.50548672488
$ws.Range("S4").Value = 0.02166461085308095
$ws.Range("T4").Value = 0.02166461085308096

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Rtn4"
$ws.Range("C5").Value = "Rtn4rl1"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 54.09018966666667
$ws.Range("H5").Value = 162.270569
$ws.Range("I5").Value = 0.2105601368412127
$ws.Range("J5").Value = 0.2105601368412127
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.620945
$ws.Range("N5").Value = 4.862835
$ws.Range("O5").Value = 0.4497939646201873
$ws.Range("P5").Value = 0.4497939646201873
$ws.Range("Q5").Value = 87.67722248923499
$ws.Range("R5").Value = 789.0950024031149
$ws.Range("S5").Value = 0.09470867874077821
$ws.Range("T5").Value = 0.09470867874077821

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Rtn4"
$ws.Range("C6").Value = "Rtn4rl1"
$ws.Range("D6").Value = "MuSCs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 54.09018966666667
$ws.Range("H6").Value = 162.270569
$ws.Range("I6").Value = 0.2105601368412127
$ws.Range("J6").Value = 0.2105601368412127
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.421245
$ws.Range("N6").Value = 4.263735
$ws.Range("O6").Value = 0.3943794658342005
$ws.Range("P6").Value = 0.3943794658342005
$ws.Range("Q6").Value = 76.87541161280166
$ws.Range("R6").Value = 691.878704515215
$ws.Range("S6").Value = 0.08304059429341362
$ws.Range("T6").Value = 0.08304059429341362

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Rtn4"
$ws.Range("C7").Value = "Rtn4rl1"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 54.09018966666667
$ws.Range("H7").Value = 162.270569
$ws.Range("I7").Value = 0.2105601368412127
$ws.Range("J7").Value = 0.2105601368412127
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.5615599999999999
$ws.Range("N7").Value = 1.68468
$ws.Range("O7").Value = 0.1558265695456122
$ws.Range("P7").Value = 0.1558265695456122
$ws.Range("Q7").Value = 30.37488690921333
$ws.Range("R7").Value = 273.37398218292
$ws.Range("S7").Value = 0.03281086380702085
$ws.Range("T7").Value = 0.03281086380702086

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Rtn4"
$ws.Range("C8").Value = "Rtn4rl1"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 101.4529346666666
$ws.Range("H8").Value = 304.358804
$ws.Range("I8").Value = 0.3949319449238378
$ws.Range("J8").Value = 0.3949319449238378
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.620945
$ws.Range("N8").Value = 4.862835
$ws.Range("O8").Value = 0.4497939646201873
$ws.Range("P8").Value = 0.4497939646201873
$ws.Range("Q8").Value = 164.44962718326
$ws.Range("R8").Value = 1480.04664464934
$ws.Range("S8").Value = 0.1776380052624544
$ws.Range("T8").Value = 0.1776380052624545

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Rtn4"
$ws.Range("C9").Value = "Rtn4rl1"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 101.4529346666666
$ws.Range("H9").Value = 304.358804
$ws.Range("I9").Value = 0.3949319449238378
$ws.Range("J9").Value = 0.3949319449238378
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.421245
$ws.Range("N9").Value = 4.263735
$ws.Range("O9").Value = 0.3943794658342005
$ws.Range("P9").Value = 0.3943794658342005
$ws.Range("Q9").Value = 144.1894761303266
$ws.Range("R9").Value = 1297.70528517294
$ws.Range("S9").Value = 0.155753049479925
$ws.Range("T9").Value = 0.155753049479925

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Rtn4"
$ws.Range("C10").Value = "Rtn4rl1"
$ws.Range("D10").Value = "Resolving-Mac"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 101.4529346666666
$ws.Range("H10").Value = 304.358804
$ws.Range("I10").Value = 0.3949319449238378
$ws.Range("J10").Value = 0.3949319449238378
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.5615599999999999
$ws.Range("N10").Value = 1.68468
$ws.Range("O10").Value = 0.1558265695456122
$ws.Range("P10").Value = 0.1558265695456122
$ws.Range("Q10").Value = 56.97190999141332
$ws.Range("R10").Value = 512.7471899227199
$ws.Range("S10").Value = 0.06154089018145829
$ws.Range("T10").Value = 0.06154089018145831

# Row 11
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "Rtn4"
$ws.Range("C11").Value = "Rtn4rl1"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 65.628919
$ws.Range("H11").Value = 196.886757
$ws.Range("I11").Value = 0.2554776429984823
$ws.Range("J11").Value = 0.2554776429984823
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 1.620945
$ws.Range("N11").Value = 4.862835
$ws.Range("O11").Value = 0.4497939646201873
$ws.Range("P11").Value = 0.4497939646201873
$ws.Range("Q11").Value = 106.380868108455
$ws.Range("R11").Value = 957.4278129760949
$ws.Range("S11").Value = 0.1149123019161082
$ws.Range("T11").Value = 0.1149123019161082

# Row 12
$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("B12").Value = "Rtn4"
$ws.Range("C12").Value = "Rtn4rl1"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 65.628919
$ws.Range("H12").Value = 196.886757
$ws.Range("I12").Value = 0.2554776429984823
$ws.Range("J12").Value = 0.2554776429984823
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 1.421245
$ws.Range("N12").Value = 4.263735
$ws.Range("O12").Value = 0.3943794658342005
$ws.Range("P12").Value = 0.3943794658342005
$ws.Range("Q12").Value = 93.27477298415499
$ws.Range("R12").Value = 839.4729568573948
$ws.Range("S12").Value = 0.100755136378322
$ws.Range("T12").Value = 0.100755136378322

# Row 13
$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("B13").Value = "Rtn4"
$ws.Range("C13").Value = "Rtn4rl1"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 65.628919
$ws.Range("H13").Value = 196.886757
$ws.Range("I13").Value = 0.2554776429984823
$ws.Range("J13").Value = 0.2554776429984823
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.5615599999999999
$ws.Range("N13").Value = 1.68468
$ws.Range("O13").Value = 0.1558265695456122
$ws.Range("P13").Value = 0.1558265695456122
$ws.Range("Q13").Value = 36.85457575363999
$ws.Range("R13").Value = 331.69118178276
$ws.Range("S13").Value = 0.03981020470405209
$ws.Range("T13").Value = 0.0398102047040521

# Remove now-obsolete rows (previously Resolving-Mac as sender; data fully superseded by rows 2-13 above)
$ws.Rows("14:17").Delete()
